$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (B2, C2 shared-string refs + F2:T2 numeric values)
# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rspo1"
$ws.Range("C2").Value = "Lgr6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03450166666666667
$ws.Range("H2").Value = 0.103505
$ws.Range("I2").Value = 0.01971168248317875
$ws.Range("J2").Value = 0.01971168248317875
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01182833333333333
$ws.Range("N2").Value = 0.035485
$ws.Range("O2").Value = 0.01157399102261438
$ws.Range("P2").Value = 0.01157399102261438
$ws.Range("Q2").Value = 0.0004080972138888889
$ws.Range("R2").Value = 0.003672874925
$ws.Range("S2").Value = 0.0002281428361009359
$ws.Range("T2").Value = 0.0002281428361009359

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rspo1"
$ws.Range("C3").Value = "Lgr6"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03450166666666667
$ws.Range("H3").Value = 0.103505
$ws.Range("I3").Value = 0.01971168248317875
$ws.Range("J3").Value = 0.01971168248317875
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.010147
$ws.Range("N3").Value = 3.030441
$ws.Range("O3").Value = 0.9884260089773856
$ws.Range("P3").Value = 0.9884260089773856
$ws.Range("Q3").Value = 0.03485175507833333
$ws.Range("R3").Value = 0.313665795705
$ws.Range("S3").Value = 0.01948353964707781
$ws.Range("T3").Value = 0.01948353964707781

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rspo1"
$ws.Range("C4").Value = "Lgr6"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.52558
$ws.Range("H4").Value = 4.57674
$ws.Range("I4").Value = 0.871602779460543
$ws.Range("J4").Value = 0.871602779460543
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01182833333333333
$ws.Range("N4").Value = 0.035485
$ws.Range("O4").Value = 0.01157399102261438
$ws.Range("P4").Value = 0.01157399102261438
$ws.Range("Q4").Value = 0.01804506876666667
$ws.Range("R4").Value = 0.1624056189
$ws.Range("S4").Value = 0.01008792274476206
$ws.Range("T4").Value = 0.01008792274476206

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rspo1"
$ws.Range("C5").Value = "Lgr6"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.52558
$ws.Range("H5").Value = 4.57674
$ws.Range("I5").Value = 0.871602779460543
$ws.Range("J5").Value = 0.871602779460543
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.010147
$ws.Range("N5").Value = 3.030441
$ws.Range("O5").Value = 0.9884260089773856
$ws.Range("P5").Value = 0.9884260089773856
$ws.Range("Q5").Value = 1.54106006026
$ws.Range("R5").Value = 13.86954054234
$ws.Range("S5").Value = 0.861514856715781
$ws.Range("T5").Value = 0.861514856715781

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Rspo1"
$ws.Range("C6").Value = "Lgr6"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.190234
$ws.Range("H6").Value = 0.570702
$ws.Range("I6").Value = 0.1086855380562782
$ws.Range("J6").Value = 0.1086855380562782
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.01182833333333333
$ws.Range("N6").Value = 0.035485
$ws.Range("O6").Value = 0.01157399102261438
$ws.Range("P6").Value = 0.01157399102261438
$ws.Range("Q6").Value = 0.002250151163333334
$ws.Range("R6").Value = 0.02025136047
$ws.Range("S6").Value = 0.001257925441751378
$ws.Range("T6").Value = 0.001257925441751377

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Rspo1"
$ws.Range("C7").Value = "Lgr6"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.190234
$ws.Range("H7").Value = 0.570702
$ws.Range("I7").Value = 0.1086855380562782
$ws.Range("J7").Value = 0.1086855380562782
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.010147
$ws.Range("N7").Value = 3.030441
$ws.Range("O7").Value = 0.9884260089773856
$ws.Range("P7").Value = 0.9884260089773856
$ws.Range("Q7").Value = 0.192164304398
$ws.Range("R7").Value = 1.729478739582
$ws.Range("S7").Value = 0.1074276126145269
$ws.Range("T7").Value = 0.1074276126145269
